$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 18.69152608107289
$ws.Range("D2").Value = 411.8440793000502
